# fix mac_ids in apis
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A now holds the mac id (device_id), column B holds the
# concatenation of "mac_id UUID" (product_uuid)
$ws.Range("A2").Value = "KEEPR90628497"
$ws.Range("B2").Value = "KEEPR90628497 EEF8EF65-AAAA-4410-B201-B6E1C4B9A486"

$ws.Range("A3").Value = "KEEPR90638498"
$ws.Range("B3").Value = "KEEPR90638498 EEF8EF65-AAAA-4410-B201-B6E1C4B9A486"

# Adjust column widths: A narrower, B wider to fit the longer value
# (input values chosen so the engine's pixel-quantized ColumnWidth
# serializes to the target OOXML widths of 15 and 31.75)
$ws.Columns.Item(1).ColumnWidth = 14.15
$ws.Columns.Item(2).ColumnWidth = 30.85

# Update the active selection to match the author's final cursor position
$ws.Range("B11").Select()
